$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column Z (day 25) values for each store row, and recompute the
# row totals in column AG (total_loja = sum of days, columns B:AF).

$ws.Range("Z2").Value = 10507.6
$ws.Range("AG2").Value = 264772.74

$ws.Range("Z3").Value = 5582
$ws.Range("AG3").Value = 120530.09

$ws.Range("Z4").Value = 1379.75
$ws.Range("AG4").Value = 77394.35000000001

$ws.Range("Z5").Value = 2711.11
$ws.Range("AG5").Value = 68371.61

$ws.Range("Z6").Value = 20180.46
$ws.Range("AG6").Value = 531068.79
